$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: extend with P1=14 and Q1=15, copying O1's formatting (bold,
# centered, bordered header style) so no new style entries are introduced.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-25: fix the I/K/M/O value pattern and append the new P/Q
# columns (value 2 in every row), matching the unstyled data-cell formatting.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2
}
